$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FAPs -> ECs)
$ws.Range("I2").Value = 0.4258685247791129
$ws.Range("J2").Value = 0.5266590341921646
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.01941266666666666
$ws.Range("N2").Value = 0.058238
$ws.Range("O2").Value = 0.01959774779686022
$ws.Range("P2").Value = 0.01978237955823044
$ws.Range("Q2").Value = 0.006706532778666666
$ws.Range("R2").Value = 0.060358795008
$ws.Range("S2").Value = 0.00834606394324197
$ws.Range("T2").Value = 0.01041856891216046

# Row 3 (FAPs -> FAPs)
$ws.Range("I3").Value = 0.4258685247791129
$ws.Range("J3").Value = 0.5266590341921646
$ws.Range("O3").Value = 0.9524028256184742
$ws.Range("P3").Value = 0.9613754796729409
$ws.Range("S3").Value = 0.4055983863415983
$ws.Range("T3").Value = 0.5063170816205801

# Row 4 (FAPs -> MuSCs)
$ws.Range("I4").Value = 0.4258685247791129
$ws.Range("J4").Value = 0.5266590341921646
$ws.Range("M4").Value = 0.027735
$ws.Range("N4").Value = 0.05547
$ws.Range("O4").Value = 0.02799942658466558
$ws.Range("P4").Value = 0.01884214076882864
$ws.Range("Q4").Value = 0.009581665919999999
$ws.Range("R4").Value = 0.05748999552
$ws.Range("S4").Value = 0.01192407449427261
$ws.Range("T4").Value = 0.009923383659424101

# Row 5 (MuSCs -> ECs)
$ws.Range("G5").Value = 0.4657455
$ws.Range("H5").Value = 0.931491
$ws.Range("I5").Value = 0.5741314752208871
$ws.Range("J5").Value = 0.4733409658078355
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01941266666666666
$ws.Range("N5").Value = 0.058238
$ws.Range("O5").Value = 0.01959774779686022
$ws.Range("P5").Value = 0.01978237955823044
$ws.Range("Q5").Value = 0.009041362142999998
$ws.Range("R5").Value = 0.05424817285799999
$ws.Range("S5").Value = 0.01125168385361825
$ws.Range("T5").Value = 0.009363810646069977

# Row 6 (MuSCs -> FAPs)
$ws.Range("G6").Value = 0.4657455
$ws.Range("H6").Value = 0.931491
$ws.Range("I6").Value = 0.5741314752208871
$ws.Range("J6").Value = 0.4733409658078355
$ws.Range("O6").Value = 0.9524028256184742
$ws.Range("P6").Value = 0.9613754796729409
$ws.Range("Q6").Value = 0.4393881859125
$ws.Range("R6").Value = 2.636329115475
$ws.Range("S6").Value = 0.5468044392768759
$ws.Range("T6").Value = 0.4550583980523609

# Row 7 (MuSCs -> MuSCs)
$ws.Range("G7").Value = 0.4657455
$ws.Range("H7").Value = 0.931491
$ws.Range("I7").Value = 0.5741314752208871
$ws.Range("J7").Value = 0.4733409658078355
$ws.Range("M7").Value = 0.027735
$ws.Range("N7").Value = 0.05547
$ws.Range("O7").Value = 0.02799942658466558
$ws.Range("P7").Value = 0.01884214076882864
$ws.Range("Q7").Value = 0.0129174514425
$ws.Range("R7").Value = 0.05166980576999999
$ws.Range("S7").Value = 0.01607535209039298
$ws.Range("T7").Value = 0.00891875710940454
